$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.710.85'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.021.50'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '510.67'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.40'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.437'
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("E10").Value = '  +1.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.368'
$ws.Range("E11").Value = '  +3.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.529.86'
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.69'
$ws.Range("E14").Value = '  +3.24%  '
$ws.Range("E15").Value = '  +5.13%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.29'
$ws.Range("E16").Value = '  +5.78%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.633.42'
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.018.67'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.95'
$ws.Range("E19").Value = '  +3.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.01'
$ws.Range("E20").Value = '  +1.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.29'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.501'
$ws.Range("E23").Value = '  +3.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.80'
$ws.Range("E24").Value = '  +3.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0925'
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.85'
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.39'
$ws.Range("E29").Value = '  +3.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.81'
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("E31").Value = '  -5.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.73'
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.78'
$ws.Range("E33").Value = '  +4.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.21'
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.91'
$ws.Range("E35").Value = '  +4.32%  '
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '24.57'
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.047.99'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.63'
$ws.Range("E40").Value = '  +1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.88'
$ws.Range("E41").Value = '  +6.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.651'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.231.57'
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.990'
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.05'
$ws.Range("E47").Value = '  +4.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0239'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("E50").Value = '  -6.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0898'
$ws.Range("E51").Value = '  +3.21%  '
